# Update column G ("K") values for rows 2-62 on the active worksheet.
# These new values represent the regenerated K (strike count) data computed
# from the underlying s_vals as described in the commit message.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newK = @{
    2  = 0
    3  = 3
    4  = 1
    5  = 2
    6  = 3
    7  = 1
    8  = 2
    9  = 0
    10 = 2
    11 = 1
    12 = 1
    13 = 0
    14 = 1
    15 = 0
    16 = 0
    17 = 1
    18 = 0
    19 = 1
    20 = 0
    21 = 2
    22 = 0
    23 = 1
    24 = 0
    25 = 1
    26 = 0
    27 = 1
    28 = 1
    29 = 1
    30 = 0
    31 = 1
    32 = 2
    33 = 1
    34 = 0
    35 = 1
    36 = 0
    37 = 0
    38 = 2
    39 = 0
    40 = 0
    41 = 1
    42 = 1
    43 = 1
    44 = 1
    45 = 0
    46 = 2
    47 = 2
    48 = 1
    49 = 2
    50 = 4
    51 = 2
    52 = 1
    53 = 3
    54 = 4
    55 = 2
    56 = 3
    57 = 3
    58 = 2
    59 = 1
    60 = 1
    61 = 2
    62 = 1
}

foreach ($row in $newK.Keys) {
    $ws.Cells.Item($row, 7).Value = $newK[$row]
}
